# Generate Report for handback
# Update the handoff/handback timestamps in the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (D) / Correspond Handback DateTime (G)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-26 06:29:55"
$wsZh.Range("D3").Value = "2016-01-26 06:29:55"
$wsZh.Range("G2").Value = "2016-01-26 06:31:32"
$wsZh.Range("G3").Value = "2016-01-26 06:31:32"

# de-de sheet: Correspond Handoff Datetime (D) / Correspond Handback DateTime (G)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-26 06:30:11"
$wsDe.Range("D3").Value = "2016-01-26 06:30:11"
$wsDe.Range("G2").Value = "2016-01-26 06:31:56"
$wsDe.Range("G3").Value = "2016-01-26 06:31:56"
